$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above current row 4, shifting rows 4-8 (and the
# placeholder rows 13/20) down by one.
$ws.Rows(4).Insert()

# The other question rows are 18pt tall (auto-fit for the 13pt Courier
# New text); match that here too.
$ws.Rows(4).RowHeight = 18

# Build the new SQL question text (question n.7: students older than 30
# years), matching the syntax-highlighted look of the other questions.
# Non-breaking spaces (U+00A0) are used as separators between the
# colored tokens, exactly like the pre-existing strings in the sheet.
$nbsp = [char]0x00A0
$newText = "SELECT" + $nbsp + "*" + $nbsp + "FROM" + $nbsp + "``students``" + $nbsp + "WHERE" + $nbsp + "TIMESTAMPDIFF(YEAR, ``date_of_birth``, CURDATE())  > 30;"

$cell = $ws.Range("B4")
$cell.Value = $newText

# Colors (hex RRGGBB -> BGR decimal expected by Font.Color)
$colGray   = 4473924    # 444444
$colPink   = 16711935   # FF00FF
$colPurple = 8913015    # 770088
$colBlue   = 11162880   # 0055AA

# "SELECT" (chars 1-6) keeps the cell's default/inherited formatting,
# so it is intentionally left untouched (no rPr in the saved file).

$cell.Characters(7, 1).Font.Color = $colGray     # " "
$cell.Characters(8, 1).Font.Color = $colPink     # "*"
$cell.Characters(9, 1).Font.Color = $colGray     # " "
$cell.Characters(10, 4).Font.Color = $colPurple  # "FROM"
$cell.Characters(14, 1).Font.Color = $colGray    # " "
$cell.Characters(15, 10).Font.Color = $colBlue   # "`students`"
$cell.Characters(25, 1).Font.Color = $colGray    # " "
$cell.Characters(26, 5).Font.Color = $colPurple  # "WHERE"
$cell.Characters(31, 55).Font.Color = $colGray   # " TIMESTAMPDIFF(...)  > 30;"

# Make the new question the active selection, as in the saved workbook.
$ws.Range("B4").Select()
